# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price updates to the Tonberry_Profits sheets
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3287
$ws.Range("J17").Value = 2759.4
$ws.Range("L17").Value = 8278.200000000001
$ws.Range("N17").Value = -8614.200000000001
$ws.Range("H18").Value = 14889.454
$ws.Range("I18").Value = 6700
$ws.Range("J18").Value = 19569.143
$ws.Range("K18").Value = 6700
$ws.Range("L18").Value = 19569.143
$ws.Range("M18").Value = -6416
$ws.Range("N18").Value = -20137.143
$ws.Range("H40").Value = 2126.889
$ws.Range("I40").Value = 1977.4286
$ws.Range("K40").Value = 1977.4286
$ws.Range("M40").Value = -1802.4286
$ws.Range("H94").Value = 2761.1428
$ws.Range("I94").Value = 2761.1428
$ws.Range("K94").Value = 2761.1428
$ws.Range("M94").Value = -2310.1428
$ws.Range("H97").Value = 921
$ws.Range("J97").Value = 925.7143
$ws.Range("L97").Value = 2777.1429
$ws.Range("N97").Value = -3769.1429
$ws.Range("H112").Value = 1803.0286
$ws.Range("J112").Value = 1821.3939
$ws.Range("L112").Value = 5464.1817
$ws.Range("N112").Value = -7680.1817
$ws.Range("H121").Value = 539.2
$ws.Range("J121").Value = 624.125
$ws.Range("L121").Value = 1872.375
$ws.Range("N121").Value = -5366.375
$ws.Range("H132").Value = 1367.6842
$ws.Range("I132").Value = 1333.8572
$ws.Range("K132").Value = 4001.5716
$ws.Range("M132").Value = -1471.5716
$ws.Range("H138").Value = 3159.532
$ws.Range("J138").Value = 2434.4666
$ws.Range("L138").Value = 7303.399800000001
$ws.Range("N138").Value = -17583.3998

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1658.3
$ws.Range("I45").Value = 1096.6666
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1096.6666
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -719.6666
$ws.Range("N45").Value = -2653
$ws.Range("H61").Value = 5371.5386
$ws.Range("I61").Value = 5274.3335
$ws.Range("K61").Value = 5274.3335
$ws.Range("M61").Value = -5062.3335
$ws.Range("H97").Value = 1385.1428
$ws.Range("I97").Value = 1299.6
$ws.Range("J97").Value = 1599
$ws.Range("K97").Value = 1299.6
$ws.Range("L97").Value = 1599
$ws.Range("M97").Value = -803.5999999999999
$ws.Range("N97").Value = -2591
$ws.Range("H132").Value = 1510.3871
$ws.Range("I132").Value = 1116.8077
$ws.Range("J132").Value = 3557
$ws.Range("K132").Value = 3350.4231
$ws.Range("L132").Value = 10671
$ws.Range("M132").Value = -820.4231
$ws.Range("N132").Value = -15731
$ws.Range("H136").Value = 5371.5386
$ws.Range("I136").Value = 5274.3335
$ws.Range("K136").Value = 15823.0005
$ws.Range("M136").Value = -13273.0005

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2139.8572
$ws.Range("I94").Value = 552.44446
$ws.Range("J94").Value = 4997.2
$ws.Range("K94").Value = 552.44446
$ws.Range("L94").Value = 4997.2
$ws.Range("M94").Value = -101.44446
$ws.Range("N94").Value = -5899.2
$ws.Range("H134").Value = 4955.8965
$ws.Range("I134").Value = 5278.84
$ws.Range("K134").Value = 15836.52
$ws.Range("M134").Value = -13301.52

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 895.3333
$ws.Range("I22").Value = 375.25
$ws.Range("K22").Value = 375.25
$ws.Range("M22").Value = -25.25
$ws.Range("H58").Value = 2417723
$ws.Range("I58").Value = 2719438.2
$ws.Range("K58").Value = 2719438.2
$ws.Range("M58").Value = -2719235.2
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -34920
$ws.Range("H132").Value = 1671.5238
$ws.Range("I132").Value = 1189.1765
$ws.Range("K132").Value = 3567.5295
$ws.Range("M132").Value = -1037.5295
$ws.Range("H134").Value = 1957.2188
$ws.Range("I134").Value = 1714.4
$ws.Range("K134").Value = 5143.200000000001
$ws.Range("M134").Value = -2608.200000000001
$ws.Range("H136").Value = 2417723
$ws.Range("I136").Value = 2719438.2
$ws.Range("K136").Value = 8158314.600000001
$ws.Range("M136").Value = -8155764.600000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 888.5
$ws.Range("I117").Value = 555
$ws.Range("J117").Value = 999.6667
$ws.Range("K117").Value = 1665
$ws.Range("L117").Value = 2999.0001
$ws.Range("M117").Value = 1777
$ws.Range("N117").Value = -9883.000100000001
$ws.Range("H131").Value = 21917.395
$ws.Range("J131").Value = 24842.896
$ws.Range("L131").Value = 74528.68799999999
$ws.Range("N131").Value = -84608.68799999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 2512.7896
$ws.Range("I102").Value = 2161
$ws.Range("K102").Value = 2161
$ws.Range("M102").Value = -539
$ws.Range("H113").Value = 951.25
$ws.Range("I113").Value = 539.25
$ws.Range("J113").Value = 1775.25
$ws.Range("K113").Value = 539.25
$ws.Range("L113").Value = 1775.25
$ws.Range("M113").Value = 1630.75
$ws.Range("N113").Value = -6115.25
$ws.Range("H122").Value = 1479.8
$ws.Range("J122").Value = 1666.6666
$ws.Range("L122").Value = 4999.9998
$ws.Range("N122").Value = -9899.9998
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1375891.1
$ws.Range("I132").Value = 1749851.5
$ws.Range("J132").Value = 4702.8335
$ws.Range("K132").Value = 5249554.5
$ws.Range("L132").Value = 14108.5005
$ws.Range("M132").Value = -5247024.5
$ws.Range("N132").Value = -19168.5005

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 13928.643
$ws.Range("I40").Value = 14260.7
$ws.Range("K40").Value = 14260.7
$ws.Range("M40").Value = -14124.7
$ws.Range("H68").Value = 1485.875
$ws.Range("I68").Value = 1485.875
$ws.Range("K68").Value = 1485.875
$ws.Range("M68").Value = -736.875
$ws.Range("H71").Value = 1485.875
$ws.Range("I71").Value = 1485.875
$ws.Range("K71").Value = 7429.375
$ws.Range("M71").Value = -3685.375
$ws.Range("H122").Value = 8728.571
$ws.Range("I122").Value = 8409.182000000001
$ws.Range("K122").Value = 25227.546
$ws.Range("M122").Value = -22777.546
$ws.Range("H132").Value = 2589.2122
$ws.Range("I132").Value = 1731.7273
$ws.Range("K132").Value = 5195.1819
$ws.Range("M132").Value = -2665.1819
$ws.Range("H136").Value = 2081.1765
$ws.Range("I136").Value = 1714.5834
$ws.Range("J136").Value = 2961
$ws.Range("K136").Value = 5143.7502
$ws.Range("L136").Value = 8883
$ws.Range("M136").Value = -2593.7502
$ws.Range("N136").Value = -13983

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 63049
$ws.Range("J42").Value = 63049
$ws.Range("L42").Value = 63049
$ws.Range("N42").Value = -63805
$ws.Range("H126").Value = 1294.9354
$ws.Range("I126").Value = 1007.2727
$ws.Range("K126").Value = 3021.8181
$ws.Range("M126").Value = -551.8181
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
